$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Schedule sheet updates
$wsSchedule.Range("E4").Value = 448.3461644999999
$wsSchedule.Range("F4").Value = 29.65252410714285
$wsSchedule.Range("E5").Value = -86.28447750000001
$wsSchedule.Range("F5").Value = -2.685480158730159

# Detailed sheet updates
$wsDetailed.Range("B39").Value = 51.95879
$wsDetailed.Range("B40").Value = 56.69657
$wsDetailed.Range("B41").Value = 58.98382
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 66.08193
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 67.15246999999999
$wsDetailed.Range("B45").Value = 69.5729
$wsDetailed.Range("B46").Value = 64.8901
$wsDetailed.Range("B48").Value = 59.75305
$wsDetailed.Range("B49").Value = 60.42201
$wsDetailed.Range("B58").Value = 57.06003
$wsDetailed.Range("B59").Value = 57.06003
$wsDetailed.Range("B60").Value = 65
$wsDetailed.Range("B61").Value = 73.20005
$wsDetailed.Range("B62").Value = 78
$wsDetailed.Range("B63").Value = 65
$wsDetailed.Range("B65").Value = 8.790279999999999
$wsDetailed.Range("B66").Value = 0.7
$wsDetailed.Range("B67").Value = 0.7
$wsDetailed.Range("B68").Value = 0.51
$wsDetailed.Range("B69").Value = -2.83936
$wsDetailed.Range("B70").Value = -6.61424
$wsDetailed.Range("B71").Value = -6.76834
$wsDetailed.Range("B72").Value = -7.92889
$wsDetailed.Range("B73").Value = -7.9236
$wsDetailed.Range("B74").Value = -7.93067
$wsDetailed.Range("B75").Value = -8.85615
$wsDetailed.Range("B76").Value = -9.99
$wsDetailed.Range("B77").Value = -8.29476
$wsDetailed.Range("B78").Value = -8.170629999999999
$wsDetailed.Range("B79").Value = -9.51909
$wsDetailed.Range("B80").Value = -8
$wsDetailed.Range("B81").Value = -6.36145
$wsDetailed.Range("B82").Value = -2.87144
$wsDetailed.Range("B83").Value = -5.51
$wsDetailed.Range("B84").Value = -0.89855
$wsDetailed.Range("B85").Value = 9.476100000000001
$wsDetailed.Range("B86").Value = 9.46405
$wsDetailed.Range("B87").Value = 33.13937
$wsDetailed.Range("B88").Value = 46.27115
$wsDetailed.Range("B89").Value = 64.8901
$wsDetailed.Range("B90").Value = 64.8901
$wsDetailed.Range("B91").Value = 58.9363
$wsDetailed.Range("B92").Value = 58.57981
$wsDetailed.Range("B93").Value = 57.6198
$wsDetailed.Range("B95").Value = 58.48911
$wsDetailed.Range("B96").Value = 61.43258
